$d = $word.ActiveDocument

# Locate the end of the "Map editor - Use bitmap or pixel" bullet text
# (right before the trailing bookmark / paragraph mark) without relying
# on hard-coded character offsets.
$findRange = $d.Content
$found = $findRange.Find.Execute("Use bitmap or pixel")
if (-not $found) {
    throw "Could not find target bullet text"
}

# Split the paragraph right after the found text (before the _GoBack
# bookmark, which stays attached to the end-of-document run) so a brand
# new list paragraph is created in its place.
$splitPoint = $d.Range($findRange.End, $findRange.End)
$splitPoint.InsertBefore([char]13)

# The newly created paragraph is now the last paragraph in the document.
$newPara = $d.Paragraphs.Last

# Insert the new bullet text at the very start of that paragraph (ahead
# of the bookmark tags that now live there).
$insertPoint = $d.Range($newPara.Range.Start, $newPara.Range.Start)
$insertPoint.InsertBefore("Implemented sprites")

# Match the Arial font used throughout the list.
$textRange = $d.Range($newPara.Range.Start, $newPara.Range.Start + 20)
$textRange.Font.Name = "Arial"

# The new bullet sits one level up from its "Map editor" sibling (ilvl 0
# vs 1); ListLevelNumber is 1-based, so 1 == w:ilvl val="0".
$newPara.Range.ListFormat.ListLevelNumber = 1
